$d = $word.ActiveDocument

# Delete the first table (the "Verantwortlich / Information / Naechste Pruefung / Typ" metadata table)
$d.Tables(1).Delete()

# Delete the title paragraph ("Dokumentationstitel") plus the empty
# paragraphs that used to surround the table, leaving just a single empty
# paragraph before the "Test Ü1" heading.
$start = $d.Paragraphs(1).Range.Start
$end = $d.Paragraphs(5).Range.End
$d.Range($start, $end).Delete()
